$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.506.78'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.12%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.790.77'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.09%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '328.55'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.59%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.08%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4399'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.45%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3759'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.19%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.72'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.45%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07633'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.99%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.141'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.15%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.66'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.39%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.005'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.12%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.266'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.31%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.461'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.25%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.791.28'
$ws.Range('D16').Style = 'Normal'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001093'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.90%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06720'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.58%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '81.81'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.88%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.002'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.01%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.58'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.63%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.246'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.93%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.565.52'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.26%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.83'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.12%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.444'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.58%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '20.62'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.95%  '

$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '153.49'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.09%  '

$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.379'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.44%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.996.40'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.16%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.314'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.27%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '130.83'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.91%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.970'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.41%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.900'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.85%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.09296'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.47%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.2251'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.55%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.20'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.32%  '

$ws.Range('B37').Value = 'TheSandbox'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6700'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.89%  '

$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06332'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.95%  '

$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02339'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.45%  '

$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.230'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.35%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.205'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.21%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.446'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.62%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.071'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.07%  '

$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.14'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.86%  '

$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.09%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6145'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.53%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.810'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.13%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '128.23'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.05%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.030'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.65%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07011'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.00%  '

$ws.Range('E51').Value = '  -1.37%  '
